# Black Scholes with setting parameters in Variables environment
#
# 1. Rename the "Input 3M" sheet to "Input_3M".
# 2. Make "Input_3M" the active sheet/tab (was "Scenario Generator").
# 3. Move the selection on "Input_3M" to cell F25.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Input 3M")
$ws.Name = "Input_3M"

$ws.Activate()
$ws.Range("F25").Select()
